$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dialogue rows (ID, Text, Speaker) to append to the "Text Lines - main"
# table. Comments column is left blank, matching every other row above.
$newRows = @(
    @("cycles_Cycles_LineTest_65J9",               "This should be fine.",          "FRED"),
    @("cycles_Cycles_LineTest_XFQW",                "So should this.",               "GEORGE"),
    @("cycles_Cycles_FancyBarkTest_RR4G",            "Fancy Bark 1",                  "FRED"),
    @("cycles_Cycles_FancyBarkTest_D4KV",            "Fancy Bark 2",                  "FRED"),
    @("cycles_Cycles_FancyBarkTest_A2I1",            "Fancy Bark 3",                  "FRED"),
    @("cycles_Cycles_FancyBarkTest_3KK1",            "Fancy Bark 4",                  "FRED"),
    @("cycles_Cycles_FancyBarkTest_FF35",            "Spinning on fancy bark 5",      "FRED"),
    @("cycles_Cycles_FancyBarkTest_23Q8",            "Spinning on fancy bark 6",      "FRED"),
    @("cycles_Cycles_StringExpressionsTest_1L9A",    "Huh.",                          "GEORGE"),
    @("cycles_Cycles_StringExpressionsTest_ZHNZ",    "Huh yourself.",                 "FRED"),
    @("cycles_Cycles_ListExpressionTest_LUCG",       "List item 1.",                  "GEORGE"),
    @("cycles_Cycles_ListExpressionTest_JXXD",       "List item 2.",                  "GEORGE")
)

$startRow = 79
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i

    # Match the row height used throughout the rest of the sheet.
    $ws.Rows.Item($r).RowHeight = 13.5

    $ws.Range("A$r").Value = $newRows[$i][0]
    $ws.Range("B$r").Value = $newRows[$i][1]
    $ws.Range("C$r").Value = $newRows[$i][2]
}

# Grow the table (and its autofilter) to cover the newly added rows.
$lo = $ws.ListObjects.Item("Table1")
$lastRow = ($startRow + $newRows.Count) - 1
$newTableRef = "A1:D$lastRow"
$lo.Resize($ws.Range($newTableRef))

Write-Host "New used range:" $ws.UsedRange.Address()
Write-Host "New table range:" $lo.Range.Address()
